# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund holdings for that quarter) between
# the existing "2020-Q4" sheet and the "总计" (grand-total) summary sheet,
# then prepends a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

$sheetQ4 = $wb.Worksheets.Item("2020-Q4")

# ---- 1. Create the new "2022-Q1" sheet right after "2020-Q4" ----------
$newSheet = $wb.Worksheets.Add($null, $sheetQ4)
$newSheet.Name = "2022-Q1"

# Re-resolve the "总计" sheet AFTER inserting the new tab so the reference
# reflects the updated sheet order (stale references to it can otherwise
# keep pointing at the wrong worksheet slot).
$sheetTotal = $wb.Worksheets.Item("总计")

# Use 总计's existing header/index-column style as the template, since the
# new sheet reuses that same visual style (bold, centered, thin border).
$headerStyle = $sheetTotal.Range("B1:D1").Style
$indexStyle = $sheetTotal.Range("A2").Style

# ---- Header row ---------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Style = $headerStyle

# ---- Data rows (columns B-G are stored as text, matching the source) ---
$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "000968"
$newSheet.Cells.Item(2, 3).Value = "广发中证养老产业指数A"
$newSheet.Cells.Item(2, 4).Value = "10.39"
$newSheet.Cells.Item(2, 5).Value = "94.08"
$newSheet.Cells.Item(2, 6).Value = "1.46"
$newSheet.Cells.Item(2, 7).Value = "0.1517"
$newSheet.Cells.Item(2, 8).Value = 8

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "002982"
$newSheet.Cells.Item(3, 3).Value = "广发中证养老产业指数C"
$newSheet.Cells.Item(3, 4).Value = "0.88"
$newSheet.Cells.Item(3, 5).Value = "94.08"
$newSheet.Cells.Item(3, 6).Value = "1.46"
$newSheet.Cells.Item(3, 7).Value = "0.0128"
$newSheet.Cells.Item(3, 8).Value = 8

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "516560"
$newSheet.Cells.Item(4, 3).Value = "华宝养老ETF"
$newSheet.Cells.Item(4, 4).Value = "0.75"
$newSheet.Cells.Item(4, 5).Value = "97.92"
$newSheet.Cells.Item(4, 6).Value = "1.52"
$newSheet.Cells.Item(4, 7).Value = "0.0114"
$newSheet.Cells.Item(4, 8).Value = 8

$newSheet.Range("A2:A4").Style = $indexStyle

# ---- 2. Prepend a "2022-Q1" summary row to "总计" ----------------------
# Existing data rows (starting at row 2) shift down by one; row indices in
# column A (0, 1, 2, ...) are renumbered afterwards to stay sequential.
$sheetTotal.Rows.Item(2).Insert()

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 3
$sheetTotal.Cells.Item(2, 4).Value = 0.18
$sheetTotal.Range("A2").Style = $indexStyle

$sheetTotal.Cells.Item(3, 1).Value = 1
$sheetTotal.Range("A3").Style = $indexStyle
